$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert new "2022-Q4" sheet right before the existing "2022-Q3" sheet
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-look-up "2022-Q3" fresh (the previously captured reference tracks
# whatever sheet now sits at that original index, not the name) so we
# copy formatting from the correct, still-existing quarter sheet.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Copy header row + column-A number style from the donor sheet so the
# new sheet matches the existing look (bold/centered header, bordered
# index column).
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Header labels
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data columns B, D..G are stored as text (numbers-as-strings) like the
# other quarter sheets, so force text format before assigning values -
# otherwise numeric-looking strings ("012072", "2.02", ...) get silently
# coerced to real numbers (and "012072" would also lose its leading
# zero). Column C (fund name) is non-numeric text already, so it is
# left alone to avoid stamping an unnecessary number-format style on it.
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "012072"
$newSheet.Range("C2").Value = "中加喜利回报一年持有期混合C"
$newSheet.Range("D2").Value = "2.02"
$newSheet.Range("E2").Value = "36.36"
$newSheet.Range("F2").Value = "1.59"
$newSheet.Range("G2").Value = "0.0321"
$newSheet.Range("H2").Value = 10

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "012071"
$newSheet.Range("C3").Value = "中加喜利回报一年持有期混合A"
$newSheet.Range("D3").Value = "1.82"
$newSheet.Range("E3").Value = "36.36"
$newSheet.Range("F3").Value = "1.59"
$newSheet.Range("G3").Value = "0.0289"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    ahead of the existing quarters, shifting the rest down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows("2").Insert()

# Re-apply the bordered/centered index-column style to the new A2 (the
# row Insert leaves it unstyled) by copying the format from A3, then
# clear the leftover header-row style that Insert copied into B2:D2.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.06

# The row-3..7 "A" index column is a fixed positional counter (0..5),
# not data that travels with the shifted row - Insert() dragged the old
# value down together with the row, so restore the correct sequence.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
